$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K) for rows 2-16, per the commit's regeneration
# of save_data using K instead of Strike#.
$values = @(3, 2, 3, 1, 1, 0, 3, 1, 4, 2, 0, 4, 2, 2, 2)

$row = 2
foreach ($v in $values) {
    $ws.Range("G$row").Value = $v
    $row++
}
